$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the formatting
# used by the existing header row (copy format from H1, the last header cell).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the values since PasteSpecial(formats) should not touch them,
# but make sure they are correct regardless of paste behavior.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-40
$i0 = @(7,6,6,5,8,4,7,8,9,6,5,8,7,6,7,6,8,8,7,9,7,8,7,4,8,8,8,7,9,6,7,8,5,5,6,5,6,4,3)
$if = @(7,6,6,5,8,6,8,8,9,7,6,8,8,6,7,6,8,8,7,9,8,9,8,6,8,9,8,7,9,7,8,8,5,5,6,5,6,4,3)

for ($idx = 0; $idx -lt $i0.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0[$idx]
    $ws.Cells.Item($row, 10).Value = $if[$idx]
}
